$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.984.34'
$ws.Range("E2").Value = '  +3.39%  '

$ws.Range("D3").Value = '3.295.65'
$ws.Range("E3").Value = '  +2.34%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '399.86'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.84%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.61'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.13%  '

$ws.Range("E7").Value = '  +5.39%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.998'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.10%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.626'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.51%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.58'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.06%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0967'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.73%  '

$ws.Range("E12").Value = '  +1.55%  '

$ws.Range("D13").Value = '3.798.96'
$ws.Range("E13").Value = '  +1.89%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.33'
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '19.10'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.49%  '

$ws.Range("D16").Value = '3.283.90'
$ws.Range("E16").Value = '  +1.98%  '

$ws.Range("E17").Value = '  -0.33%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.01'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.79%  '

$ws.Range("D19").Value = '57.726.54'
$ws.Range("E19").Value = '  +3.24%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.34'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.29%  '

$ws.Range("E21").Value = '  +5.48%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '13.01'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.50%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '300.17'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.48%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.64'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.04%  '

$ws.Range("E25").Value = '  -0.49%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '28.31'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.75%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.94'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.71%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.40'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.77%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.42'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.38%  '

$ws.Range("E30").Value = '  -1.29%  '

$ws.Range("E31").Value = '  -0.44%  '

$ws.Range("E32").Value = '  +2.58%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.31'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.78%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '40.73'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +12.72%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0500'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.02%  '

$ws.Range("E36").Value = '  +0.79%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.86'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.92%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.18'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.58%  '

$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.51'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.44%  '

$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.04%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '138.68'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.27%  '

$ws.Range("E42").Value = '  +2.18%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.89'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.27%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.284'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.23%  '

$ws.Range("B45").Value = 'Celestia'
$ws.Range("C45").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.94'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.39%  '

$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.92'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.66%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.39'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.90%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.23'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.52%  '

$ws.Range("D49").Value = '2.166.26'
$ws.Range("E49").Value = '  +2.25%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.47'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.13%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.92'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -11.08%  '
